# Remove the "Bắc Kinh (Trung Quốc)", "Thượng Hải (Trung Quốc)" and
# "Tokyo (Nhật Bản)" columns from the summary table.
#
# Original column layout (row 1 headers):
#   A Năm | B Jakarta | C Bangkok | D Singapore | E Manila | F Kuala Lumpur |
#   G Hong Kong | H Seoul | I Taipei | J Bắc Kinh | K Thượng Hải |
#   L Thâm Quyến | M Tokyo | N Osaka | O TP. Hồ Chí Minh
#
# After deleting J (Bắc Kinh), K (Thượng Hải) and M (Tokyo) the remaining
# columns shift left so that:
#   J Thâm Quyến | K Osaka | L TP. Hồ Chí Minh

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from right to left so earlier deletions don't shift the letters
# of the columns still to be removed.
$ws.Columns("M").Delete()
$ws.Columns("K").Delete()
$ws.Columns("J").Delete()

# Keep the previously-selected cell pointing at the bottom-right corner of
# the (now smaller) data range, same as Excel would after deleting columns.
$ws.Range("L6").Select() | Out-Null
